function Find-ParagraphIndex($doc, $targetText, $occurrence) {
    $seen = 0
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
        if ($t -eq $targetText) {
            $seen++
            if ($seen -eq $occurrence) {
                return $i
            }
        }
    }
    return -1
}

function Insert-OpenXmlParagraphs($range, $bodyXml) {
    $pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`n" + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
    $pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkgHeader + $bodyXml + $pkgFooter)
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Before the (first) "$ git status" paragraph that follows "... first
#    commit", insert two new paragraphs:
#       $ git add effsite.R
#       $ git status
# ---------------------------------------------------------------------
$statusIdx = Find-ParagraphIndex $d '$ git status' 1
if ($statusIdx -eq -1) { throw "could not locate '$ git status' paragraph" }

$statusPara = $d.Paragraphs.Item($statusIdx)
$statusPara.Range.InsertParagraphBefore()

# The freshly-created (still empty) paragraph is now at $statusIdx;
# replace its content (and trailing mark) with the two real paragraphs.
$newPara = $d.Paragraphs.Item($statusIdx)
$insertRange = $d.Range($newPara.Range.Start, $newPara.Range.End)

$addParagraph = '<w:p><w:r><w:t xml:space="preserve">$ </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>git</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> add </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>effsite.R</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '</w:p>'

$statusParagraph = '<w:p><w:r><w:t xml:space="preserve">$ </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>git</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> status</w:t></w:r>' + `
    '</w:p>'

$twoParagraphs = $addParagraph + $statusParagraph
Insert-OpenXmlParagraphs $insertRange $twoParagraphs

# ---------------------------------------------------------------------
# 2) Remove the <w:lastRenderedPageBreak/> from the (second) "$ git log"
#    paragraph.
# ---------------------------------------------------------------------
$logIdx = Find-ParagraphIndex $d '$ git log' 2
if ($logIdx -eq -1) { throw "could not locate second '$ git log' paragraph" }

$logPara = $d.Paragraphs.Item($logIdx)
$logRange = $d.Range($logPara.Range.Start, $logPara.Range.End)

$logParagraph = '<w:p w:rsidR="006458C2" w:rsidRDefault="006458C2" w:rsidP="006458C2">' + `
    '<w:r><w:t xml:space="preserve">$ </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>git</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> log</w:t></w:r>' + `
    '</w:p>'

Insert-OpenXmlParagraphs $logRange $logParagraph

# ---------------------------------------------------------------------
# 3) Delete the whole "$ git remote add origin master" paragraph.
# ---------------------------------------------------------------------
$remoteIdx = Find-ParagraphIndex $d '$ git remote add origin master' 1
if ($remoteIdx -eq -1) { throw "could not locate '$ git remote add origin master' paragraph" }

$remotePara = $d.Paragraphs.Item($remoteIdx)
$remotePara.Range.Delete()

Write-Output "done"
